# Generate Report for Handoff
#
# Re-running the localization status report refreshed the "Latest Handoff
# Datetime" (column D) for the row belonging to the
# "1dd41960-bd9a-4302-ac86-61dfb1d8037e" file on both the zh-cn and de-de
# sheets, while everything else in the report stayed the same.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-10 12:45:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-10 12:45:20"
